$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.379.44'
$ws.Range("E2").Value = '  +1.30%  '

$ws.Range("D3").Value = '1.821.95'
$ws.Range("E3").Value = '  -0.10%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.82%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4480'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.43%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3758'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.30%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07499'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.23%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8871'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.11%  '

$ws.Range("E11").Value = '  +1.65%  '

$ws.Range("D12").Value = '1.827.67'
$ws.Range("E12").Value = '  +0.06%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.757'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.70%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.95'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.412'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07096'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.37%  '

$ws.Range("E17").Value = '  -0.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008801'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("E20").Value = '  +1.88%  '

$ws.Range("D21").Value = '27.389.60'
$ws.Range("E21").Value = '  +0.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.271'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.93'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.39%  '

$ws.Range("D24").Value = '2.058.22'
$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.960'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.375'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.78%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.59'
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.391'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.19%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '118.01'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.95%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08850'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.84%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7884'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.202'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.513'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.11%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.908'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9997'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.115'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.75%  '

$ws.Range("E38").Value = '  +2.56%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05326'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.84%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.366'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.92%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5322'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.45%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1725'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.74%  '

$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.856'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.59%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.314'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +19.56%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.737'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.39%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5100'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.67'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.45%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.703'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.83%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.78'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9998'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.09%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06375'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.59%  '
